$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.231825238986011
$ws.Range("K2").Value = 0.0928686071010423
$ws.Range("L2").Value = 0.384326285472465
$ws.Range("N2").Value = 0.1767713663003

$ws.Range("B3").Value = 0.221756139772863
$ws.Range("K3").Value = 0.202637183919994
$ws.Range("L3").Value = 0.113737513135779
$ws.Range("N3").Value = 0.257273414011137

$ws.Range("B4").Value = 0.318130319930611
$ws.Range("K4").Value = 0.451904268645967
$ws.Range("L4").Value = 0.227706886529004
$ws.Range("N4").Value = 0.338155063004691

$ws.Range("B5").Value = 0.172200105515648
$ws.Range("E5").Value = 0.16615846782681
$ws.Range("K5").Value = 0.19166544338377
$ws.Range("L5").Value = 0.194374624175739
$ws.Range("N5").Value = 0.175591451815085
